# pl_mw.xlsx / Sheet1 — rerun of the power-flow results for the 380 kV case.
# Rows 2-25 hold one result row per timestep (col A = 0..23); only the
# B,D,E,F,G,H,I,J,N result columns change, the rest (A,C,K,L,M,O) stay put.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newResults = @(
    @{ Row = 2; "B" = 0.2159311196965064; "D" = 0.1058422808887993; "E" = 0.13028179933346; "F" = 2.878478165862106; "G" = 2.368769996223534; "H" = 1.696868061313182; "I" = 0.3757315010511491; "J" = 0.2151472402880472; "N" = 1.538965158848526 }
    @{ Row = 3; "B" = 0.1913949238938244; "D" = 0.1057099897463161; "E" = 0.129376248143025; "F" = 2.750918434283136; "G" = 2.222640816841505; "H" = 1.636692919714051; "I" = 0.3656238581545637; "J" = 0.208773148833032; "N" = 1.440536212657179 }
    @{ Row = 4; "B" = 0.1762869034274956; "D" = 0.1056316089112226; "E" = 0.1288919000434596; "F" = 2.674709327376632; "G" = 2.13468666062019; "H" = 1.600961959854857; "I" = 0.3596127975951617; "J" = 0.205012059613054; "N" = 1.380350891855841 }
    @{ Row = 5; "B" = 0.170120147310584; "D" = 0.1056003683078792; "E" = 0.1287124574646867; "F" = 2.64417696227224; "G" = 2.099280920985308; "H" = 1.586703403558658; "I" = 0.3572128864356259; "J" = 0.2035173224339673; "N" = 1.3558914054328 }
    @{ Row = 6; "B" = 0.1690955695415823; "D" = 0.1055952226743999; "E" = 0.128683741779934; "F" = 2.639138474748648; "G" = 2.093427924020915; "H" = 1.584353915396832; "I" = 0.3568173988915788; "J" = 0.2032714010193857; "N" = 1.351834048475013 }
    @{ Row = 7; "B" = 0.1762037765289932; "D" = 0.1056311847747633; "E" = 0.1288894075201767; "F" = 2.674295448035082; "G" = 2.134207411306051; "H" = 1.600768445849013; "I" = 0.3595802297553021; "J" = 0.2049917480354679; "N" = 1.380020748618364 }
    @{ Row = 8; "B" = 0.2074802657996884; "D" = 0.1057960665953246; "E" = 0.129954637554043; "F" = 2.834052162487836; "G" = 2.318011996819791; "H" = 1.675864740954125; "I" = 0.3722062484302811; "J" = 0.212917504022812; "N" = 1.504976964631908 }
    @{ Row = 9; "B" = 0.2684501888316504; "D" = 0.1061427228871796; "E" = 0.1326164825695173; "F" = 3.164481609023056; "G" = 2.692914156986546; "H" = 1.832970800449971; "I" = 0.3984893586954854; "J" = 0.2296922894002904; "N" = 1.751860150726543 }
    @{ Row = 10; "B" = 0.3129948669854343; "D" = 0.106412707435446; "E" = 0.1349279334950673; "F" = 3.418265018242977; "G" = 2.977790851650184; "H" = 1.954674957977772; "I" = 0.4186982167615696; "J" = 0.2427996710749767; "N" = 1.934196595620278 }
    @{ Row = 11; "B" = 0.3331993617556179; "D" = 0.106539096441006; "E" = 0.1360582219594733; "F" = 3.536236937461979; "G" = 3.10957700108986; "H" = 2.011467405619328; "I" = 0.4280806917633271; "J" = 0.2489397106595561; "N" = 2.01731744366225 }
    @{ Row = 12; "B" = 0.3408412057246153; "D" = 0.1065874905544781; "E" = 0.1364976832432454; "F" = 3.581282987014646; "G" = 3.15980762891769; "H" = 2.033183669984112; "I" = 0.4316602591565868; "J" = 0.2512908683668513; "N" = 2.04881500450972 }
    @{ Row = 13; "B" = 0.3391958157679369; "D" = 0.1065770439752107; "E" = 0.1364025264829394; "F" = 3.571564784942183; "G" = 3.148974891615353; "H" = 2.02849725803037; "I" = 0.4308881605544457; "J" = 0.2507833367248367; "N" = 2.042030543333738 }
    @{ Row = 14; "B" = 0.3338282482070554; "D" = 0.1065430670494472; "E" = 0.1360941466472845; "F" = 3.539935378736061; "G" = 3.113702908760729; "H" = 2.013249775342672; "I" = 0.4283746547068858; "J" = 0.2491326164057455; "N" = 2.019908357050213 }
    @{ Row = 15; "B" = 0.3305392456624645; "D" = 0.1065223252554794; "E" = 0.1359067491117649; "F" = 3.520610273165971; "G" = 3.092140582485115; "H" = 2.003937769423828; "I" = 0.4268385096919118; "J" = 0.2481249124802929; "N" = 2.006360570313319 }
    @{ Row = 16; "B" = 0.3116732098857256; "D" = 0.1064045211602469; "E" = 0.1348556618023871; "F" = 3.410606892642278; "G" = 2.969223366258007; "H" = 1.950992589575719; "I" = 0.4180888088536463; "J" = 0.2424020215006948; "N" = 1.928767661017872 }
    @{ Row = 17; "B" = 0.3000838993125114; "D" = 0.1063331807884911; "E" = 0.1342311160739662; "F" = 3.343776744307405; "G" = 2.894387060466613; "H" = 1.918881687400301; "I" = 0.4127692216154841; "J" = 0.2389370388622893; "N" = 1.881209144293678 }
    @{ Row = 18; "B" = 0.2934125106850445; "D" = 0.1062924826548013; "E" = 0.1338793036318329; "F" = 3.305574844941731; "G" = 2.851549089412629; "H" = 1.900546521073124; "I" = 0.4097273953839675; "J" = 0.2369607392022033; "N" = 1.853871441568089 }
    @{ Row = 19; "B" = 0.2911527646054139; "D" = 0.1062787599065302; "E" = 0.1337614551878872; "F" = 3.292680760349469; "G" = 2.837079934219787; "H" = 1.89436145048802; "I" = 0.4087005711172651; "J" = 0.2362944442791957; "N" = 1.844618345498816 }
    @{ Row = 20; "B" = 0.3013181776852036; "D" = 0.1063407402842653; "E" = 0.1342968323355507; "F" = 3.350866326720222; "G" = 2.90233212842918; "H" = 1.922286027038183; "I" = 0.4133336567285824; "J" = 0.2393041633429505; "N" = 1.886270128092036 }
    @{ Row = 21; "B" = 0.3354050867557703; "D" = 0.1065530322609014; "E" = 0.136184413718567; "F" = 3.549215510681506; "G" = 3.124054210283646; "H" = 2.017722587296021; "I" = 0.429112214418673; "J" = 0.249616761358638; "N" = 2.026405625114307 }
    @{ Row = 22; "B" = 0.3576291578725375; "D" = 0.1066948957068323; "E" = 0.1374848268051494; "F" = 3.681024872415463; "G" = 3.270868098814447; "H" = 2.081323391753187; "I" = 0.4395793436413982; "J" = 0.2565087896800833; "N" = 2.118115898853603 }
    @{ Row = 23; "B" = 0.3457728835221303; "D" = 0.1066188884862704; "E" = 0.1367846235793699; "F" = 3.610473322205905; "G" = 3.192332840114375; "H" = 2.047264520717476; "I" = 0.4339788619097504; "J" = 0.2528162789516131; "N" = 2.069158297353852 }
    @{ Row = 24; "B" = 0.3007601869816199; "D" = 0.10633732165204; "E" = 0.1342670994488167; "F" = 3.347660443095805; "G" = 2.898739583346128; "H" = 1.920746533548197; "I" = 0.4130784242512391; "J" = 0.2391381373017794; "N" = 1.883982043575145 }
    @{ Row = 25; "B" = 0.2519980490933733; "D" = 0.1060463370010964; "E" = 0.1318343853988502; "F" = 3.073198553213103; "G" = 2.589881758271019; "H" = 1.789387422476921; "I" = 0.3912191043998661; "J" = 0.2250191059075775; "N" = 1.684892220676261 }
)

foreach ($entry in $newResults) {
    foreach ($col in "B", "D", "E", "F", "G", "H", "I", "J", "N") {
        $addr = "$col$($entry.Row)"
        $ws.Range($addr).Value = $entry[$col]
    }
}
